$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4.2
$ws.Range("I2").Value = 1.94
$ws.Range("V2").Value = 2.06
$ws.Range("W2").Value = 1.21
$ws.Range("I3").Value = 3.3
$ws.Range("J3").Value = 2.72
$ws.Range("K3").Value = 2.84
$ws.Range("L3").Value = 1.6
$ws.Range("X3").Value = 6.6
$ws.Range("Y3").Value = 8.199999999999999
$ws.Range("AB3").Value = 7.6
$ws.Range("AC3").Value = 6.8
$ws.Range("AM3").Value = 300
$ws.Range("AO3").Value = 1000
$ws.Range("F4").Value = 1.7
$ws.Range("H4").Value = 5.1
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 3.8
$ws.Range("O4").Value = 1.29
$ws.Range("P4").Value = 1.99
$ws.Range("Q4").Value = 1.84
$ws.Range("S4").Value = 3.2
$ws.Range("T4").Value = 1.82
$ws.Range("U4").Value = 1.98
$ws.Range("V4").Value = 1.2
$ws.Range("W4").Value = 2.26
$ws.Range("Y4").Value = 19
$ws.Range("Z4").Value = 46
$ws.Range("AB4").Value = 9
$ws.Range("AD4").Value = 21
$ws.Range("AF4").Value = 11.5
$ws.Range("AG4").Value = 10.5
$ws.Range("AJ4").Value = 19
$ws.Range("AK4").Value = 19.5
$ws.Range("F5").Value = 1.68
$ws.Range("G5").Value = 1.76
$ws.Range("K5").Value = 4.1
$ws.Range("P5").Value = 1.81
$ws.Range("W5").Value = 2.3
$ws.Range("G6").Value = 1.27
$ws.Range("H6").Value = 14.5
$ws.Range("I6").Value = 17.5
$ws.Range("J6").Value = 6.8
$ws.Range("K6").Value = 8.199999999999999
$ws.Range("L6").Value = 1.22
$ws.Range("N6").Value = 5.7
$ws.Range("O6").Value = 1.16
$ws.Range("P6").Value = 2.6
$ws.Range("Q6").Value = 1.5
$ws.Range("R6").Value = 1.64
$ws.Range("S6").Value = 2.28
$ws.Range("T6").Value = 2.08
$ws.Range("U6").Value = 1.74
$ws.Range("X6").Value = 36
$ws.Range("AB6").Value = 11
$ws.Range("AC6").Value = 17.5
$ws.Range("AE6").Value = 300
$ws.Range("AF6").Value = 8.800000000000001
$ws.Range("AG6").Value = 12
$ws.Range("AI6").Value = 210
$ws.Range("AJ6").Value = 9.800000000000001
$ws.Range("AM6").Value = 210
$ws.Range("J7").Value = 7.6
$ws.Range("AC7").Value = 970
$ws.Range("AO7").Value = 3.7
$ws.Range("L8").Value = 1.56
$ws.Range("Z8").Value = 24
$ws.Range("AD8").Value = 18.5
$ws.Range("O9").Value = 1.49
$ws.Range("S9").Value = 4.3
$ws.Range("T9").Value = 2.22
$ws.Range("V9").Value = 1.17
$ws.Range("Y9").Value = 19
$ws.Range("AJ9").Value = 1000
$ws.Range("AN9").Value = 1000
$ws.Range("F10").Value = 1.71
$ws.Range("G10").Value = 1.84
$ws.Range("K10").Value = 4
$ws.Range("W10").Value = 2.2
$ws.Range("Z10").Value = 46
$ws.Range("AN10").Value = 12.5
$ws.Range("L11").Value = 1.58
$ws.Range("F12").Value = 2.24
$ws.Range("G12").Value = 2.32
$ws.Range("I12").Value = 3.8
$ws.Range("N12").Value = 3.5
$ws.Range("P12").Value = 1.84
$ws.Range("Q12").Value = 2.1
$ws.Range("T12").Value = 1.84
$ws.Range("U12").Value = 2.1
$ws.Range("V12").Value = 1.36
$ws.Range("W12").Value = 1.75
$ws.Range("AA12").Value = 75
$ws.Range("AC12").Value = 7.6
$ws.Range("AD12").Value = 16
$ws.Range("AE12").Value = 50
$ws.Range("AG12").Value = 11.5
$ws.Range("AH12").Value = 20
$ws.Range("AI12").Value = 60
$ws.Range("AM12").Value = 130
$ws.Range("G13").Value = 2.4
$ws.Range("I13").Value = 4.4
$ws.Range("J13").Value = 3.05
$ws.Range("V13").Value = 1.29
$ws.Range("W13").Value = 1.71
$ws.Range("G14").Value = 1.7
$ws.Range("I14").Value = 6.8
$ws.Range("N14").Value = 3.55
$ws.Range("Q14").Value = 2.04
$ws.Range("R14").Value = 1.34
$ws.Range("T14").Value = 2
$ws.Range("AD14").Value = 24
$ws.Range("F15").Value = 1.46
$ws.Range("G15").Value = 1.5
$ws.Range("H15").Value = 8.6
$ws.Range("J15").Value = 4.4
$ws.Range("N15").Value = 3.4
$ws.Range("P15").Value = 1.82
$ws.Range("S15").Value = 3.85
$ws.Range("U15").Value = 1.69
$ws.Range("W15").Value = 3
$ws.Range("AE15").Value = 230
$ws.Range("AH15").Value = 38
$ws.Range("AM15").Value = 280
$ws.Range("AO15").Value = 410
$ws.Range("F16").Value = 1.47
$ws.Range("I16").Value = 11
$ws.Range("J16").Value = 4
$ws.Range("L16").Value = 1.4
$ws.Range("V16").Value = 1.1
$ws.Range("Y16").Value = 24
$ws.Range("AD16").Value = 42
$ws.Range("AJ16").Value = 14
$ws.Range("AK16").Value = 22
$ws.Range("AN16").Value = 12.5
